# smoke test cases done
# Adds 15 new smoke-test rows (33-47) to Sheet1, covering the
# TourToPK / Content Manager / Hotel Manager / Tour Operator dashboards.
#
# NOTE: the ".Value = ..." assignments below for columns A and B are
# intentionally issued in this exact (non-row-sequential) order so the
# workbook's shared-string table is rebuilt in the same order the
# original author typed the cells in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$xlCenter = -4108

$ws.Range("B33").Value = "Able to navigate all links"
$ws.Range("B35").Value = "able to view and delete plans"
$ws.Range("B34").Value = "able to view and delete places"
$ws.Range("A33").Value = "TourToPK Dashboard"
$ws.Range("B36").Value = "able to view and delete hotels"
$ws.Range("B37").Value = "able to view and delete packages"
$ws.Range("B39").Value = "able to create partners accounts"
$ws.Range("A40").Value = "Content Manager`ndashboad"
$ws.Range("B40").Value = "able to navigate all links "
$ws.Range("B41").Value = "able to create, update, delete places"
$ws.Range("B42").Value = "abe to create,update,delete plans"
$ws.Range("B38").Value = "able to view and delete all tourists"
$ws.Range("A43").Value = "Hotel Manger `ndahsboard"
$ws.Range("B43").Value = "able to naviage all links"
$ws.Range("B44").Value = "able to create, update, delete hotels"
$ws.Range("B45").Value = "able to view and delete booking"
$ws.Range("A46").Value = "Tour Operator `ndashboard"
$ws.Range("B46").Value = "able to naviage all links"
$ws.Range("B47").Value = "able to create, update, delete packages"

# --- Result column (reuses the existing "PASS" shared string) -------------
$ws.Range("C33").Value = "PASS"
$ws.Range("C34").Value = "PASS"
$ws.Range("C35").Value = "PASS"
$ws.Range("C36").Value = "PASS"
$ws.Range("C37").Value = "PASS"
$ws.Range("C38").Value = "PASS"
$ws.Range("C39").Value = "PASS"
$ws.Range("C40").Value = "PASS"
$ws.Range("C41").Value = "PASS"
$ws.Range("C42").Value = "PASS"
$ws.Range("C43").Value = "PASS"
$ws.Range("C44").Value = "PASS"
$ws.Range("C45").Value = "PASS"
$ws.Range("C46").Value = "PASS"
$ws.Range("C47").Value = "PASS"

# --- Formatting: section-header rows get centered/wrapped text style ------
$ws.Range("A33").HorizontalAlignment = $xlCenter

$ws.Range("A40").WrapText = $true
$ws.Rows.Item(40).RowHeight = 30

$ws.Range("A43").WrapText = $true
$ws.Rows.Item(43).RowHeight = 30

$ws.Range("A46").WrapText = $true
$ws.Rows.Item(46).RowHeight = 30

# --- Data validation range grows along with the new rows ------------------
$ws.Range("C2:C47").Validation.Delete()
$ws.Range("C2:C47").Validation.Add(3, 1, 1, "PASS, FAIL")
$ws.Range("C2:C47").Validation.InCellDropdown = $true
$ws.Range("C2:C47").Validation.IgnoreBlank = $true

# --- Move the on-screen selection to match the author's saved view --------
$ws.Range("D51").Select()
